$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New style for K1: bold font (same as header font), no border, centered/top aligned
$ws.Range("B1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Borders.LineStyle = -4142
$excel.CutCopyMode = $false

# Column widths
$ws.Columns.Item(1).ColumnWidth = 3.5924479166666665
$ws.Columns.Item(12).ColumnWidth = 11.022135416666666

# View state
$ws.Range("L5").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
